$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 7
$ws.Range("H7").Value = 22224.5
$ws.Range("J7").Value = 22224.5
$ws.Range("L7").Value = 22224.5
$ws.Range("N7").Value = -22448.5

# ALC row 14
$ws.Range("H14").Value = 22224.5
$ws.Range("J14").Value = 22224.5
$ws.Range("L14").Value = 22224.5
$ws.Range("N14").Value = -22606.5

# ALC row 17
$ws.Range("H17").Value = 521.6491
$ws.Range("J17").Value = 361.91837
$ws.Range("L17").Value = 1085.75511
$ws.Range("N17").Value = -1421.75511

# ALC row 64
$ws.Range("H64").Value = 3825
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3052
$ws.Range("N64").Value = -4496

# ALC row 67
$ws.Range("H67").Value = 3825
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2442
$ws.Range("N67").Value = -5716

# ALC row 81
$ws.Range("H81").Value = 31060
$ws.Range("J81").Value = 31060
$ws.Range("L81").Value = 31060
$ws.Range("N81").Value = -33056

# ALC row 84
$ws.Range("H84").Value = 31060
$ws.Range("J84").Value = 31060
$ws.Range("L84").Value = 93180
$ws.Range("N84").Value = -103164

# ALC row 116
$ws.Range("H116").Value = 440106.9
$ws.Range("I116").Value = 1001626
$ws.Range("J116").Value = 8169.154
$ws.Range("K116").Value = 1001626
$ws.Range("L116").Value = 8169.154
$ws.Range("M116").Value = -998184
$ws.Range("N116").Value = -15053.154

# ALC row 131
$ws.Range("H131").Value = 2670.3572
$ws.Range("I131").Value = 599.1667
$ws.Range("J131").Value = 4223.75
$ws.Range("K131").Value = 1797.5001
$ws.Range("L131").Value = 12671.25
$ws.Range("M131").Value = 3242.4999
$ws.Range("N131").Value = -22751.25

# ALC row 141
$ws.Range("H141").Value = 112982.555
$ws.Range("I141").Value = 155398.53
$ws.Range("J141").Value = 2701
$ws.Range("K141").Value = 466195.59
$ws.Range("L141").Value = 8103
$ws.Range("M141").Value = -461015.59
$ws.Range("N141").Value = -18463

$ws = $wb.Worksheets.Item("ARM")
# ARM row 30
$ws.Range("H30").Value = 5500
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# ARM row 110
$ws.Range("H110").Value = 1288.4736
$ws.Range("I110").Value = 1280.3636
$ws.Range("J110").Value = 1299.625
$ws.Range("K110").Value = 1280.3636
$ws.Range("L110").Value = 1299.625
$ws.Range("M110").Value = 764.6364000000001
$ws.Range("N110").Value = -5389.625

# ARM row 132
$ws.Range("H132").Value = 2352.739
$ws.Range("I132").Value = 989.6923
$ws.Range("J132").Value = 4124.7
$ws.Range("K132").Value = 2969.0769
$ws.Range("L132").Value = 12374.1
$ws.Range("M132").Value = -439.0769
$ws.Range("N132").Value = -17434.1

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 4832541
$ws.Range("I105").Value = 5209938
$ws.Range("J105").Value = 1861.4
$ws.Range("K105").Value = 5209938
$ws.Range("L105").Value = 1861.4
$ws.Range("M105").Value = -5208191
$ws.Range("N105").Value = -5355.4

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 15873815
$ws.Range("I16").Value = 22222940
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 22222940
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -22222653
$ws.Range("N16").Value = -1574

# CRP row 74
$ws.Range("H74").Value = 32472
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 32472
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").Value = 32472
$ws.Range("N74").Value = -34220

# CRP row 77
$ws.Range("H77").Value = 32472
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 32472
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").Value = 97416
$ws.Range("N77").Value = -106152

# CRP row 100
$ws.Range("H100").Value = 68000
$ws.Range("J100").Value = 68000
$ws.Range("L100").Value = 68000
$ws.Range("N100").Value = -70164

# CRP row 113
$ws.Range("H113").Value = 15873815
$ws.Range("I113").Value = 22222940
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 22222940
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -22220770
$ws.Range("N113").Value = -5340

# CRP row 132
$ws.Range("H132").Value = 5463.625
$ws.Range("I132").Value = 5350.6
$ws.Range("J132").Value = 5652
$ws.Range("K132").Value = 16051.8
$ws.Range("L132").Value = 16956
$ws.Range("M132").Value = -13521.8
$ws.Range("N132").Value = -22016

# CRP row 137
$ws.Range("H137").Value = 46047.145
$ws.Range("J137").Value = 46047.145
$ws.Range("L137").Value = 46047.145
$ws.Range("N137").Value = -56247.145

$ws = $wb.Worksheets.Item("CUL")
# CUL row 25
$ws.Range("H25").Value = 2109.75
$ws.Range("I25").Value = 2199
$ws.Range("J25").Value = 2080
$ws.Range("K25").Value = 6597
$ws.Range("L25").Value = 6240
$ws.Range("M25").Value = -6428
$ws.Range("N25").Value = -6578

# CUL row 30
$ws.Range("H30").Value = 2109.75
$ws.Range("I30").Value = 2199
$ws.Range("J30").Value = 2080
$ws.Range("K30").Value = 6597
$ws.Range("L30").Value = 6240
$ws.Range("M30").Value = -6495
$ws.Range("N30").Value = -6444

# CUL row 50
$ws.Range("H50").Value = 390.36365
$ws.Range("I50").Value = 265.5
$ws.Range("J50").Value = 540.2
$ws.Range("K50").Value = 796.5
$ws.Range("L50").Value = 1620.6
$ws.Range("M50").Value = -315.5
$ws.Range("N50").Value = -2582.6

# CUL row 53
$ws.Range("H53").Value = 390.36365
$ws.Range("I53").Value = 265.5
$ws.Range("J53").Value = 540.2
$ws.Range("K53").Value = 796.5
$ws.Range("L53").Value = 1620.6
$ws.Range("M53").Value = -315.5
$ws.Range("N53").Value = -2582.6

# CUL row 107
$ws.Range("H107").Value = 125492.375
$ws.Range("I107").Value = 482.5
$ws.Range("J107").Value = 250502.25
$ws.Range("K107").Value = 1447.5
$ws.Range("L107").Value = 751506.75
$ws.Range("M107").Value = 472.5
$ws.Range("N107").Value = -755346.75

# CUL row 131
$ws.Range("H131").Value = 773.165
$ws.Range("I131").Value = 344.66666
$ws.Range("J131").Value = 801.4176
$ws.Range("K131").Value = 1033.99998
$ws.Range("L131").Value = 2404.2528
$ws.Range("M131").Value = 4006.00002
$ws.Range("N131").Value = -12484.2528

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Range("H70").Value = 6524.875
$ws.Range("I70").Value = 5843
$ws.Range("K70").Value = 5843
$ws.Range("M70").Value = -5573

# GSM row 73
$ws.Range("H73").Value = 6524.875
$ws.Range("I73").Value = 5843
$ws.Range("K73").Value = 5843
$ws.Range("M73").Value = -4907

$ws = $wb.Worksheets.Item("LTW")
# LTW row 30
$ws.Range("H30").Value = 5333
$ws.Range("I30").Value = 666.6667
$ws.Range("J30").Value = 9999.333000000001
$ws.Range("K30").Value = 666.6667
$ws.Range("L30").Value = 9999.333000000001
$ws.Range("M30").Value = -558.6667
$ws.Range("N30").Value = -10215.333

# LTW row 81
$ws.Range("H81").Value = 51856.855
$ws.Range("J81").Value = 51856.855
$ws.Range("L81").Value = 51856.855
$ws.Range("N81").Value = -53852.855

# LTW row 84
$ws.Range("H84").Value = 51856.855
$ws.Range("J84").Value = 51856.855
$ws.Range("L84").Value = 155570.565
$ws.Range("N84").Value = -165554.565

# LTW row 132
$ws.Range("H132").Value = 12782.774
$ws.Range("I132").Value = 19163.143
$ws.Range("J132").Value = 7528.353
$ws.Range("K132").Value = 57489.429
$ws.Range("L132").Value = 22585.059
$ws.Range("M132").Value = -54959.429
$ws.Range("N132").Value = -27645.059

# LTW row 135
$ws.Range("H135").Value = 49333.332
$ws.Range("J135").Value = 49333.332
$ws.Range("L135").Value = 49333.332
$ws.Range("N135").Value = -59473.332

$ws = $wb.Worksheets.Item("WVR")
# WVR row 29
$ws.Range("H29").Value = 36673.668
$ws.Range("I29").Value = 20005
$ws.Range("K29").Value = 20005
$ws.Range("M29").Value = -19715

# WVR row 80
$ws.Range("H80").Value = 44142.715
$ws.Range("J80").Value = 44142.715
$ws.Range("L80").Value = 44142.715
$ws.Range("N80").Value = -46138.715

# WVR row 83
$ws.Range("H83").Value = 44142.715
$ws.Range("J83").Value = 44142.715
$ws.Range("L83").Value = 132428.145
$ws.Range("N83").Value = -142412.145

# WVR row 132
$ws.Range("H132").Value = 6668802
$ws.Range("I132").Value = 1269.2433
$ws.Range("K132").Value = 3807.7299
$ws.Range("M132").Value = -1277.7299
